$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "PWM_Voltage" sheet (was tab-selected with B6) ---
$pwm = $wb.Worksheets.Item("PWM_Voltage")
$pwm.Range("A1:B1").Select() | Out-Null

# --- Add the new "Sheet1" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "Sheet1"

# Column D is a bit wider to fit the "Duty Cycle" header
$new.Columns.Item(4).ColumnWidth = 15.71

# --- Header row ---
$new.Cells.Item(1,1).Value = "PWM"
$new.Cells.Item(1,2).Value = "Voltage"
$new.Cells.Item(1,3).Value = "RPM"
$new.Cells.Item(1,4).Value = "Duty Cycle"

# --- Data rows: PWM, Voltage, RPM, DriverDuty ---
$data = @(
  @(5,4.96,0,1),
  @(8,4.74,0,6),
  @(9,4.6500000000000004,3,7),
  @(10,4.58,4.5,9),
  @(11,4.5199999999999996,6,10),
  @(12,4.46,7.5,11),
  @(13,4.41,8,12),
  @(14,4.3600000000000003,9.5,13),
  @(14,4.32,11,14),
  @(16,4.4400000000000004,11.5,15),
  @(17,4.2300000000000004,13.5,16),
  @(19,4.18,15.8,17),
  @(20,4.13,16,18),
  @(21,4.09,18,19),
  @(23,4.0199999999999996,20.5,20),
  @(24,3.99,20.5,21),
  @(26,3.92,23,22),
  @(28,3.86,23,23),
  @(29,3.82,25.5,24),
  @(31,3.76,27.8,25),
  @(33,3.7,28,26),
  @(34,3.67,30,27),
  @(35,3.64,32,28),
  @(38,3.56,33,29),
  @(39,3.53,35,30),
  @(40,3.5,37,31),
  @(42,3.44,37.5,32),
  @(44,3.39,40,33),
  @(47,3.3,42,35),
  @(51,3.2,45,37),
  @(54,3.1,46.5,38),
  @(56,3.06,48,39),
  @(58,3,50,40)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $new.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Put the selection/active cell on D2, matching the freshly-entered log
$new.Range("D2").Select() | Out-Null
